# "Updated Code for JBL Login Page"
#
# On the MyAccountLogin sheet the "Status" column (C) is duplicated into a
# new column D:
#   - C1 ("Status" header) is copied so D1 also reads "Status"
#   - C2:C4 (the Fail/Pass/Fail result values) are moved over to D2:D4,
#     so column C keeps only its header and column D holds the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MyAccountLogin")

# Move the result values out of C2:C4 into D2:D4.
$ws.Range("C2:C4").Cut() | Out-Null
$ws.Range("D2").Select() | Out-Null
$ws.Paste() | Out-Null

# Duplicate the "Status" header from C1 into the new D1 cell.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").Select() | Out-Null
$ws.Paste() | Out-Null

# Give the newly introduced header cell its own style entry.
$ws.Range("D1").HorizontalAlignment = 1

# Leave the sheet selection on the (still present) original header cell.
$ws.Range("C1").Select() | Out-Null
